$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 14.61399645465972
$ws.Range("C2").Value = 7.967680387001448
$ws.Range("D2").Value = 7.919997271376202
$ws.Range("E2").Value = 10.35691967998533
$ws.Range("F2").Value = 45.35077606733246
$ws.Range("H2").Value = 7.344005520526261
$ws.Range("M2").Value = 15.7655862068672
$ws.Range("B3").Value = 14.20314885177105
$ws.Range("C3").Value = 7.550689470233192
$ws.Range("D3").Value = 7.757603831157134
$ws.Range("E3").Value = 10.19212663783078
$ws.Range("F3").Value = 43.91823071386122
$ws.Range("H3").Value = 7.344005520526261
$ws.Range("M3").Value = 15.56563962829284
$ws.Range("B4").Value = 13.95270309439353
$ws.Range("C4").Value = 7.287496905188032
$ws.Range("D4").Value = 7.656282840507525
$ws.Range("E4").Value = 10.09202247293478
$ws.Range("F4").Value = 43.02161795793949
$ws.Range("H4").Value = 7.344005520526261
$ws.Range("M4").Value = 15.44798760829663
$ws.Range("B5").Value = 13.85130646400733
$ws.Range("C5").Value = 7.178642268386019
$ws.Range("D5").Value = 7.614621340415656
$ws.Range("E5").Value = 10.05154339788976
$ws.Range("F5").Value = 42.65244284402661
$ws.Range("H5").Value = 7.344005520526261
$ws.Range("M5").Value = 15.40138909386646
$ws.Range("B6").Value = 13.83451573891846
$ws.Range("C6").Value = 7.160476716383392
$ws.Range("D6").Value = 7.607681964732717
$ws.Range("E6").Value = 10.04484209197392
$ws.Range("F6").Value = 42.59092699773035
$ws.Range("H6").Value = 7.344005520526261
$ws.Range("M6").Value = 15.39373438137266
$ws.Range("B7").Value = 13.95133266733586
$ws.Range("C7").Value = 7.286035053115839
$ws.Range("D7").Value = 7.655722442757189
$ws.Range("E7").Value = 10.09147522978959
$ws.Range("F7").Value = 43.01665384794162
$ws.Range("H7").Value = 7.344005520526261
$ws.Range("M7").Value = 15.44735364028421
$ws.Range("B8").Value = 14.47208025210431
$ws.Range("C8").Value = 7.825497506956331
$ws.Range("D8").Value = 7.864357575617988
$ws.Range("E8").Value = 10.29990080036263
$ws.Range("F8").Value = 44.86064629784912
$ws.Range("H8").Value = 7.344005520526261
$ws.Range("M8").Value = 15.69561931726203
$ws.Range("B9").Value = 15.49905265961376
$ws.Range("C9").Value = 8.835613099371662
$ws.Range("D9").Value = 8.259337938601663
$ws.Range("E9").Value = 10.71529438401306
$ws.Range("F9").Value = 48.32227871434085
$ws.Range("H9").Value = 7.344005520526261
$ws.Range("M9").Value = 16.2203090866832
$ws.Range("B10").Value = 16.24626537978608
$ws.Range("C10").Value = 9.649016824642093
$ws.Range("D10").Value = 8.539166244345161
$ws.Range("E10").Value = 11.02202505737257
$ws.Range("F10").Value = 50.74812194904791
$ws.Range("H10").Value = 7.344005520526261
$ws.Range("M10").Value = 16.62502324401619
$ws.Range("B11").Value = 16.58255383145433
$ws.Range("C11").Value = 9.998869313032365
$ws.Range("D11").Value = 8.663879833171627
$ws.Range("E11").Value = 11.16137461780933
$ws.Range("F11").Value = 51.82202434326928
$ws.Range("H11").Value = 7.344005520526261
$ws.Range("M11").Value = 16.81246874501943
$ws.Range("B12").Value = 16.70922376080227
$ws.Range("C12").Value = 10.12846254508553
$ws.Range("D12").Value = 8.710708272926771
$ws.Range("E12").Value = 11.21407371468418
$ws.Range("F12").Value = 52.22412176610536
$ws.Range("H12").Value = 7.344005520526261
$ws.Range("M12").Value = 16.88385865377028
$ws.Range("B13").Value = 16.68197546769214
$ws.Range("C13").Value = 10.10068057884713
$ws.Range("D13").Value = 8.70064103264896
$ws.Range("E13").Value = 11.20272778640344
$ws.Range("F13").Value = 52.13773048566713
$ws.Range("H13").Value = 7.344005520526261
$ws.Range("M13").Value = 16.8684665334768
$ws.Range("B14").Value = 16.59298941054018
$ws.Range("C14").Value = 10.00958883940216
$ws.Range("D14").Value = 8.667740557827921
$ws.Range("E14").Value = 11.1657118402309
$ws.Range("F14").Value = 51.85519828490875
$ws.Range("H14").Value = 7.344005520526261
$ws.Range("M14").Value = 16.81833426045331
$ws.Range("B15").Value = 16.5383906011996
$ws.Range("C15").Value = 9.953416751492764
$ws.Range("D15").Value = 8.647535488905353
$ws.Range("E15").Value = 11.14302819007464
$ws.Range("F15").Value = 51.68153600634554
$ws.Range("H15").Value = 7.344005520526261
$ws.Range("M15").Value = 16.78767785682347
$ws.Range("B16").Value = 16.22420174549669
$ws.Range("C16").Value = 9.62574828375333
$ws.Range("D16").Value = 8.530961650494389
$ws.Range("E16").Value = 11.01291081942547
$ws.Range("F16").Value = 50.67731682762441
$ws.Range("H16").Value = 7.344005520526261
$ws.Range("M16").Value = 16.61283428138545
$ws.Range("B17").Value = 16.03041517521332
$ws.Range("C17").Value = 9.419576630585892
$ws.Range("D17").Value = 8.458766635722705
$ws.Range("E17").Value = 10.93300984206427
$ws.Range("F17").Value = 50.05344450202853
$ws.Range("H17").Value = 7.344005520526261
$ws.Range("M17").Value = 16.50637739833062
$ws.Range("B18").Value = 15.91862167939059
$ws.Range("C18").Value = 9.299092943341904
$ws.Range("D18").Value = 8.417000454920482
$ws.Range("E18").Value = 10.88703811801311
$ws.Range("F18").Value = 49.69183718955453
$ws.Range("H18").Value = 7.344005520526261
$ws.Range("M18").Value = 16.44546478204066
$ws.Range("B19").Value = 15.88071812535215
$ws.Range("C19").Value = 9.257972703874715
$ws.Range("D19").Value = 8.402818502962484
$ws.Range("E19").Value = 10.87147167780874
$ws.Range("F19").Value = 49.5689368743487
$ws.Range("H19").Value = 7.344005520526261
$ws.Range("M19").Value = 16.4248977258309
$ws.Range("B20").Value = 16.05107966891176
$ws.Range("C20").Value = 9.441720432927857
$ws.Range("D20").Value = 8.466477112935875
$ws.Range("E20").Value = 10.94151728562519
$ws.Range("F20").Value = 50.12014603780012
$ws.Range("H20").Value = 7.344005520526261
$ws.Range("M20").Value = 16.51767746645629
$ws.Range("B21").Value = 16.61914625995398
$ws.Range("C21").Value = 10.03642297790319
$ws.Range("D21").Value = 8.677415227164563
$ws.Range("E21").Value = 11.1765865453445
$ws.Range("F21").Value = 51.93831103557747
$ws.Range("H21").Value = 7.344005520526261
$ws.Range("M21").Value = 16.83304881029597
$ws.Range("B22").Value = 16.98641298227479
$ws.Range("C22").Value = 10.40826654424311
$ws.Range("D22").Value = 8.812945143224651
$ws.Range("E22").Value = 11.32979350999246
$ws.Range("F22").Value = 53.09985357738945
$ws.Range("H22").Value = 7.344005520526261
$ws.Range("M22").Value = 17.04150974892201
$ws.Range("B23").Value = 16.7908090758491
$ws.Range("C23").Value = 10.21134222308644
$ws.Range("D23").Value = 8.74083185520845
$ws.Range("E23").Value = 11.24807677111076
$ws.Range("F23").Value = 52.48245505621409
$ws.Range("H23").Value = 7.344005520526261
$ws.Range("M23").Value = 16.93005894861462
$ws.Range("B24").Value = 16.0417384324587
$ws.Range("C24").Value = 9.431715303496652
$ws.Range("D24").Value = 8.462992014939614
$ws.Range("E24").Value = 10.93767117795156
$ws.Range("F24").Value = 50.08999939115615
$ws.Range("H24").Value = 7.344005520526261
$ws.Range("M24").Value = 16.51256779422906
$ws.Range("B25").Value = 15.22177272251131
$ws.Range("C25").Value = 8.558233268203608
$ws.Range("D25").Value = 8.154185318953399
$ws.Range("E25").Value = 10.60246756448858
$ws.Range("F25").Value = 47.40489543236775
$ws.Range("H25").Value = 7.344005520526261
$ws.Range("M25").Value = 16.0747060106832
